# Applies the "Intermediate changes to frameside, C#, and tester classes"
# commit: updates a handful of Pin # values on the "Frameside" sheet
# (two of which switch from numeric pin numbers to the textual "A2"/"A5"
# analog-pin labels), and refreshes the sheet's view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Frameside")
$ws.Activate()

# Row 6 - Drawing Linac Feedback Pin: Pin # 2 -> "A2" (analog pin label)
$ws.Cells.Item(6, 4).Value = "A2"

# Row 8 - Drawing Linac Motor Direction: Pin # 1 -> 5
$ws.Cells.Item(8, 4).Value = 5

# Row 9 - Aiming Linac Feedback Pin: Pin # 5 -> "A5" (analog pin label)
$ws.Cells.Item(9, 4).Value = "A5"

# Row 14 - Pin # 20 -> 21
$ws.Cells.Item(14, 4).Value = 21

# Row 15 - Pin # 21 -> 20
$ws.Cells.Item(15, 4).Value = 20

# Refresh the view: scroll back to the top and move the selection to F14
$ws.Range("A1").Select()
$ws.Range("F14").Select()
